$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values on row 2
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("I2").Value = "Cupón"
$ws.Range("T2").Value = 21004016

# Reset the view: scroll back to A1 and select A2
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
